# This workbook keeps a weekly price log sorted with the newest entries at
# the top of the data block. A new week's worth of data (2 rows) needs to be
# inserted right after the header/previous-latest block, i.e. before the
# current row 389, pushing everything else (rows 389..419) down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 389 (existing row 389 and everything
# below shifts down by 2).
$ws.Rows.Item(389).Insert()
$ws.Rows.Item(389).Insert()

# --- New row 389 ---
$ws.Cells.Item(389, 1).Value = 4
$ws.Cells.Item(389, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(389, 3).Value = "Los Lagos"
$ws.Cells.Item(389, 4).Value = 45013
$ws.Cells.Item(389, 5).Value = 10
$ws.Cells.Item(389, 6).Value = "Fruta"
$ws.Cells.Item(389, 7).Value = 100101
$ws.Cells.Item(389, 8).Value = "Berries"
$ws.Cells.Item(389, 9).Value = 100112025
$ws.Cells.Item(389, 10).Value = "Frutilla"
$ws.Cells.Item(389, 11).Value = "Sin especificar"
$ws.Cells.Item(389, 12).Value = "Primera"
$ws.Cells.Item(389, 13).Value = 600
$ws.Cells.Item(389, 14).Value = 9500
$ws.Cells.Item(389, 15).Value = 10000
$ws.Cells.Item(389, 16).Value = 9750
$ws.Cells.Item(389, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(389, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(389, 19).Value = 1393
$ws.Cells.Item(389, 20).Value = 7

# --- New row 390 ---
$ws.Cells.Item(390, 1).Value = 4
$ws.Cells.Item(390, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(390, 3).Value = "Los Lagos"
$ws.Cells.Item(390, 4).Value = 45013
$ws.Cells.Item(390, 5).Value = 10
$ws.Cells.Item(390, 6).Value = "Fruta"
$ws.Cells.Item(390, 7).Value = 100101
$ws.Cells.Item(390, 8).Value = "Berries"
$ws.Cells.Item(390, 9).Value = 100112025
$ws.Cells.Item(390, 10).Value = "Frutilla"
$ws.Cells.Item(390, 11).Value = "Sin especificar"
$ws.Cells.Item(390, 12).Value = "Segunda"
$ws.Cells.Item(390, 13).Value = 300
$ws.Cells.Item(390, 14).Value = 8000
$ws.Cells.Item(390, 15).Value = 8000
$ws.Cells.Item(390, 16).Value = 8000
$ws.Cells.Item(390, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(390, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(390, 19).Value = 1143
$ws.Cells.Item(390, 20).Value = 7

# Make sure column D on the two new rows keeps the date number format used
# throughout the rest of the column (same style as the row above/below).
$ws.Range("D389:D390").NumberFormat = $ws.Range("D391").NumberFormat
